# Applies the "6.0.0 / Alvearie Team / Jurisdiction" metadata refresh to the
# employee-count StructureDefinition workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" -----------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Publication date refresh
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$meta.Range("B9").Value = "Alvearie Team"

# The old duplicated "Contact" row becomes "Jurisdiction"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Remove the leftover duplicate "Contact" row (old row 11); everything
# below shifts up by one, restoring the Description/Purpose/... block.
$meta.Rows.Item(11).Delete()

# --- Sheet 2: "Elements" ------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The root Extension row's Short/Definition now reflect the actual
# extension instead of the generic placeholder text.
$elements.Range("K2").Value = "Employee Count"
$elements.Range("L2").Value = "Flag of the employees in the eligibility record. Each employee is identified with a value of 1. All others are given values of 0."
